$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# Row structure changes
# -----------------------------------------------------------------
# Insert a new row at row 5 (pushes the "Source" row from 5->6, and
# the "Note" row from 6->7). The new blank row 5 inherits formatting
# from row 4 above it.
$ws.Rows.Item(5).Insert() | Out-Null

# Delete the old "Note" row, which is now at row 7 (it is removed
# entirely in the new layout).
$ws.Rows.Item(7).Delete() | Out-Null

# -----------------------------------------------------------------
# Row 1: title
# -----------------------------------------------------------------
$ws.Range("A1").Value = "The number of persons with disabilities registered in the Unified database of targeted social assistance program in Sighnaghi Municipality"
$titleRange = $ws.Range("A1:I1")
$titleRange.Merge() | Out-Null
$titleRange.HorizontalAlignment = -4108
$titleRange.VerticalAlignment = -4108
$titleRange.WrapText = $true
$titleRange.Font.Name = "Arial"
$titleRange.Font.Size = 11
$titleRange.Font.Bold = $true

# -----------------------------------------------------------------
# Row 4: "family with disabilities Persons" data row
# -----------------------------------------------------------------
$ws.Range("A4").Value = "family with disabilities Persons "
$a4 = $ws.Range("A4")
$a4.Font.Name = "Arial"
$a4.Font.Size = 10
$a4.HorizontalAlignment = -4131
$a4.VerticalAlignment = -4108
$a4.WrapText = $true
$a4.Borders.Item(8).LineStyle = 1
$a4.Borders.Item(8).Weight = 2
$a4.Borders.Item(9).LineStyle = -4142

$rowVals4 = @(537, 525, 475, 488, 484, 475, 494, 514)
for ($i = 0; $i -lt $rowVals4.Length; $i++) {
    $col = 2 + $i
    $c = $ws.Cells.Item(4, $col)
    $c.Value = $rowVals4[$i]
    $c.NumberFormat = "#\ ##0"
    $c.Borders.Item(8).LineStyle = -4142
    $c.Borders.Item(9).LineStyle = -4142
}

# -----------------------------------------------------------------
# Row 5: "disabilities Persons" data row
# -----------------------------------------------------------------
$ws.Range("A5").Value = "disabilities Persons "
$a5 = $ws.Range("A5")
$a5.Font.Name = "Arial"
$a5.Font.Size = 10
$a5.HorizontalAlignment = -4131
$a5.VerticalAlignment = -4108
$a5.WrapText = $true
$a5.Borders.Item(8).LineStyle = -4142
$a5.Borders.Item(9).LineStyle = 1
$a5.Borders.Item(9).Weight = 2

$rowVals5 = @(575, 568, 516, 531, 528, 518, 545, 566)
for ($i = 0; $i -lt $rowVals5.Length; $i++) {
    $col = 2 + $i
    $c = $ws.Cells.Item(5, $col)
    $c.Value = $rowVals5[$i]
    $c.NumberFormat = "#\ ##0"
    $c.Borders.Item(8).LineStyle = -4142
    if ($col -eq 9) {
        $c.Borders.Item(9).LineStyle = 1
        $c.Borders.Item(9).Weight = 2
    } else {
        $c.Borders.Item(9).LineStyle = -4142
    }
}

# -----------------------------------------------------------------
# Row 6: "Source: ..." row (now at row 6 after the insert above)
# -----------------------------------------------------------------
$a6 = $ws.Range("A6")
$a6.Borders.Item(8).LineStyle = -4142

# -----------------------------------------------------------------
# Column width
# -----------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 20

# -----------------------------------------------------------------
# Row heights
# -----------------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 51
$ws.Rows.Item(2).RowHeight = 14.5
$ws.Rows.Item(3).RowHeight = 18.75
$ws.Rows.Item(4).RowHeight = 24.75
$ws.Rows.Item(5).RowHeight = 21
$ws.Rows.Item(6).RowHeight = 27.75

Write-Host "done"
